# Generate Report for Handback
#
# Adds "latest target file" hyperlink + handback datetime + error detail
# for the ce6864df-... row (row 8) on both the "zh-cn" and "de-de" sheets,
# because the handback for that file is stale (not the latest commit).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7aae8785f4ec8da8bda97993153c7b5d3a64d22b/e2e/ce6864df-67db-4dae-ba11-cfa742e034a2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98bd4905b2280fa0459ab53afdf2b003ba3dc922/e2e/ce6864df-67db-4dae-ba11-cfa742e034a2.md."

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98bd4905b2280fa0459ab53afdf2b003ba3dc922/e2e/ce6864df-67db-4dae-ba11-cfa742e034a2.md"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Columns.Item(16).ColumnWidth = 40

$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $hyperlinkTarget, "", "", "ce6864df-67db-4dae-ba11-cfa742e034a2.md")
$wsZh.Range("I8").Value = "ce6864df-67db-4dae-ba11-cfa742e034a2.md"

$wsZh.Range("J8").Value = "ce6864df-67db-4dae-ba11-cfa742e034a2.39ed58b684b045a329a6c4b87e10502d6f3e418c.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-28 20:43:02"
$wsZh.Range("P8").Value = $errorDetail

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = 40

$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $hyperlinkTarget, "", "", "ce6864df-67db-4dae-ba11-cfa742e034a2.md")
$wsDe.Range("I8").Value = "ce6864df-67db-4dae-ba11-cfa742e034a2.md"

$wsDe.Range("J8").Value = "ce6864df-67db-4dae-ba11-cfa742e034a2.39ed58b684b045a329a6c4b87e10502d6f3e418c.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-28 20:43:10"
$wsDe.Range("P8").Value = $errorDetail
